# Auto-generated Excel COM-interop script to apply Spriggan_Profits.xlsx data update
# (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 547.9
$ws.Range("I2").Value = 148.45454
$ws.Range("J2").Value = 1036.1111
$ws.Range("K2").Value = 148.45454
$ws.Range("L2").Value = 1036.1111
$ws.Range("M2").Value = -35.45454000000001
$ws.Range("N2").Value = -1262.1111

$ws.Range("H9").Value = 6665.778
$ws.Range("I9").Value = 7452.2666
$ws.Range("J9").Value = 2733.3333
$ws.Range("K9").Value = 7452.2666
$ws.Range("L9").Value = 2733.3333
$ws.Range("M9").Value = -7283.2666
$ws.Range("N9").Value = -3071.3333

$ws.Range("H19").Value = 1335.0588
$ws.Range("I19").Value = 1342.091
$ws.Range("J19").Value = 1322.1666
$ws.Range("K19").Value = 1342.091
$ws.Range("L19").Value = 1322.1666
$ws.Range("M19").Value = -1167.091
$ws.Range("N19").Value = -1672.1666

$ws.Range("H21").Value = 4000
$ws.Range("I21").Value = 4000
$ws.Range("K21").Value = 4000
$ws.Range("M21").Value = -3532

$ws.Range("H23").Value = 4000
$ws.Range("I23").Value = 4000
$ws.Range("K23").Value = 4000
$ws.Range("M23").Value = -3766

$ws.Range("H29").Value = 2215.7856
$ws.Range("I29").Value = 238
$ws.Range("J29").Value = 3314.5557
$ws.Range("K29").Value = 714
$ws.Range("L29").Value = 9943.667099999999
$ws.Range("M29").Value = -433
$ws.Range("N29").Value = -10505.6671

$ws.Range("H38").Value = 1213.7778
$ws.Range("I38").Value = 1213.7778
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 3641.3334
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -3269.3334
$ws.Range("N38").ClearContents()

$ws.Range("H43").Value = 6558.625
$ws.Range("I43").Value = 6558.625
$ws.Range("K43").Value = 6558.625
$ws.Range("M43").Value = -6489.625

$ws.Range("H58").Value = 248.55556
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

$ws.Range("H76").Value = 3380
$ws.Range("I76").Value = 2818.375
$ws.Range("J76").Value = 4128.8335
$ws.Range("K76").Value = 2818.375
$ws.Range("L76").Value = 4128.8335
$ws.Range("M76").Value = -2503.375
$ws.Range("N76").Value = -4758.8335

$ws.Range("H79").Value = 3380
$ws.Range("I79").Value = 2818.375
$ws.Range("J79").Value = 4128.8335
$ws.Range("K79").Value = 2818.375
$ws.Range("L79").Value = 4128.8335
$ws.Range("M79").Value = -1726.375
$ws.Range("N79").Value = -6312.8335

$ws.Range("H112").Value = 146121.58
$ws.Range("J112").Value = 86877.336
$ws.Range("L112").Value = 260632.008
$ws.Range("N112").Value = -262848.008

$ws.Range("H135").Value = 23810098
$ws.Range("I135").Value = 23810098
$ws.Range("K135").Value = 214290882
$ws.Range("M135").Value = -214288347

$ws.Range("H137").Value = 2120.9268
$ws.Range("I137").Value = 1782.5333
$ws.Range("K137").Value = 5347.5999
$ws.Range("M137").Value = -2797.5999

$ws.Range("H138").Value = 2229.0266
$ws.Range("I138").Value = 1158.4445
$ws.Range("J138").Value = 3217.2563
$ws.Range("K138").Value = 3475.3335
$ws.Range("L138").Value = 9651.768899999999
$ws.Range("M138").Value = 1664.6665
$ws.Range("N138").Value = -19931.7689


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1636314.5
$ws.Range("I2").Value = 2102418.8
$ws.Range("J2").Value = 4949.5
$ws.Range("K2").Value = 2102418.8
$ws.Range("L2").Value = 4949.5
$ws.Range("M2").Value = -2102305.8
$ws.Range("N2").Value = -5175.5

$ws.Range("H32").Value = 2099.5957
$ws.Range("I32").Value = 1910.5405
$ws.Range("K32").Value = 1910.5405
$ws.Range("M32").Value = -1623.5405

$ws.Range("H61").Value = 35715316
$ws.Range("I61").Value = 43479150
$ws.Range("J61").Value = 1659
$ws.Range("K61").Value = 43479150
$ws.Range("L61").Value = 1659
$ws.Range("M61").Value = -43478938
$ws.Range("N61").Value = -2083

$ws.Range("H74").Value = 66675596
$ws.Range("I74").Value = 76931530
$ws.Range("J74").Value = 12000
$ws.Range("K74").Value = 76931530
$ws.Range("L74").Value = 12000
$ws.Range("M74").Value = -76930656
$ws.Range("N74").Value = -13748

$ws.Range("H77").Value = 66675596
$ws.Range("I77").Value = 76931530
$ws.Range("J77").Value = 12000
$ws.Range("K77").Value = 384657650
$ws.Range("L77").Value = 60000
$ws.Range("M77").Value = -384653282
$ws.Range("N77").Value = -68736

$ws.Range("H97").Value = 268.7619
$ws.Range("I97").Value = 296.8889
$ws.Range("J97").Value = 100
$ws.Range("K97").Value = 296.8889
$ws.Range("L97").Value = 100
$ws.Range("M97").Value = 199.1111
$ws.Range("N97").Value = -1092

$ws.Range("H116").Value = 1636314.5
$ws.Range("I116").Value = 2102418.8
$ws.Range("J116").Value = 4949.5
$ws.Range("K116").Value = 2102418.8
$ws.Range("L116").Value = 4949.5
$ws.Range("M116").Value = -2100124.8
$ws.Range("N116").Value = -9537.5

$ws.Range("H132").Value = 2224700
$ws.Range("I132").Value = 2633971.8
$ws.Range("J132").Value = 2939.7144
$ws.Range("K132").Value = 7901915.399999999
$ws.Range("L132").Value = 8819.143199999999
$ws.Range("M132").Value = -7899385.399999999
$ws.Range("N132").Value = -13879.1432

$ws.Range("H136").Value = 35715316
$ws.Range("I136").Value = 43479150
$ws.Range("J136").Value = 1659
$ws.Range("K136").Value = 130437450
$ws.Range("L136").Value = 4977
$ws.Range("M136").Value = -130434900
$ws.Range("N136").Value = -10077


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1636314.5
$ws.Range("I3").Value = 2102418.8
$ws.Range("J3").Value = 4949.5
$ws.Range("K3").Value = 2102418.8
$ws.Range("L3").Value = 4949.5
$ws.Range("M3").Value = -2102304.8
$ws.Range("N3").Value = -5177.5

$ws.Range("H26").Value = 24736.666
$ws.Range("J26").Value = 39000
$ws.Range("L26").Value = 39000
$ws.Range("N26").Value = -39584

$ws.Range("H74").Value = 16138.5
$ws.Range("J74").Value = 16138.5
$ws.Range("L74").Value = 16138.5
$ws.Range("N74").Value = -18010.5

$ws.Range("H77").Value = 16138.5
$ws.Range("J77").Value = 16138.5
$ws.Range("L77").Value = 48415.5
$ws.Range("N77").Value = -57775.5

$ws.Range("H94").Value = 13751
$ws.Range("I94").Value = 13896.728
$ws.Range("J94").Value = 12949.5
$ws.Range("K94").Value = 13896.728
$ws.Range("L94").Value = 12949.5
$ws.Range("M94").Value = -13445.728
$ws.Range("N94").Value = -13851.5

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 20005256
$ws.Range("I58").Value = 29418616
$ws.Range("K58").Value = 29418616
$ws.Range("M58").Value = -29418413

$ws.Range("H69").Value = 16000
$ws.Range("I69").Value = 16000
$ws.Range("K69").Value = 16000
$ws.Range("M69").Value = -15251

$ws.Range("H72").Value = 16000
$ws.Range("I72").Value = 16000
$ws.Range("K72").Value = 48000
$ws.Range("M72").Value = -44256

$ws.Range("H136").Value = 20005256
$ws.Range("I136").Value = 29418616
$ws.Range("K136").Value = 88255848
$ws.Range("M136").Value = -88253298


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 119418.88
$ws.Range("I11").Value = 119418.88
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 358256.64
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -358116.64
$ws.Range("N11").ClearContents()

$ws.Range("H132").Value = 1443.5555
$ws.Range("I132").Value = 1141.8572
$ws.Range("J132").Value = 2499.5
$ws.Range("K132").Value = 10276.7148
$ws.Range("L132").Value = 22495.5
$ws.Range("M132").Value = -7746.7148
$ws.Range("N132").Value = -27555.5


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 3809.2
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 5015.3335
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 5015.3335
$ws.Range("M46").Value = -1844
$ws.Range("N46").Value = -5327.3335

$ws.Range("H70").Value = 4277.25
$ws.Range("I70").Value = 4033.3333
$ws.Range("J70").Value = 5009
$ws.Range("K70").Value = 4033.3333
$ws.Range("L70").Value = 5009
$ws.Range("M70").Value = -3763.3333
$ws.Range("N70").Value = -5549

$ws.Range("H73").Value = 4277.25
$ws.Range("I73").Value = 4033.3333
$ws.Range("J73").Value = 5009
$ws.Range("K73").Value = 4033.3333
$ws.Range("L73").Value = 5009
$ws.Range("M73").Value = -3097.3333
$ws.Range("N73").Value = -6881

$ws.Range("H126").Value = 20000
$ws.Range("I126").Value = 20000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 60000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -57530
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 20835476
$ws.Range("I132").Value = 25002332
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 75006996
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -75004466
$ws.Range("N132").Value = -8660


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2710.138
$ws.Range("I16").Value = 1211.0834
$ws.Range("K16").Value = 1211.0834
$ws.Range("M16").Value = -1041.0834

$ws.Range("H43").Value = 525000
$ws.Range("J43").Value = 1000000
$ws.Range("L43").Value = 1000000
$ws.Range("N43").Value = -1000386

$ws.Range("H100").Value = 16636132
$ws.Range("I100").Value = 22179842
$ws.Range("K100").Value = 22179842
$ws.Range("M100").Value = -22179301

$ws.Range("H122").Value = 4949.75
$ws.Range("I122").Value = 4949.75
$ws.Range("K122").Value = 14849.25
$ws.Range("M122").Value = -12399.25

$ws.Range("H132").Value = 14795462
$ws.Range("I132").Value = 15781493
$ws.Range("J132").Value = 4999.6665
$ws.Range("K132").Value = 47344479
$ws.Range("L132").Value = 14998.9995
$ws.Range("M132").Value = -47341949
$ws.Range("N132").Value = -20058.9995


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 69000
$ws.Range("J120").Value = 69000
$ws.Range("L120").Value = 69000
$ws.Range("N120").Value = -78676

